$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before the current last 5 rows (98-102), shifting the
# existing rows 98-102 down to become rows 103-107 with all their original
# data/formatting intact.
$ws.Rows("98:102").Insert()

$data = @(
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44610, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Primera", 210, 10000, 10000, 10000, "`$/caja 15 kilos granel", "Región Metropolitana", 667, 15),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44610, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Segunda", 280, 8000, 8000, 8000, "`$/caja 15 kilos granel", "Región Metropolitana", 533, 15),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44610, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Tercera", 170, 5000, 5000, 5000, "`$/caja 15 kilos granel", "Región Metropolitana", 333, 15),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44610, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Larry Ann", "Primera", 85, 10000, 10000, 10000, "`$/caja 15 kilos granel", "Región de O'Higgins", 667, 15),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44610, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Larry Ann", "Segunda", 120, 8000, 8000, 8000, "`$/caja 15 kilos granel", "Región de O'Higgins", 533, 15)
)

$startRow = 98
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($col = 0; $col -lt $row.Count; $col++) {
        $ws.Cells.Item($r, $col + 1).Value = $row[$col]
    }
}
